# The tail of the document currently looks like:
#   ... <w:p/>                     (empty paragraph - rsidR 00DB1031)
#       <w:p/>                     (empty paragraph - rsidR 001C56AE)
#       <w:p>                      (paragraph containing the text "MVC")
#         <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
#         <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>MVC</w:t></w:r>
#         <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#         <w:bookmarkEnd w:id="0"/>
#       </w:p>
#
# It needs to become:
#   ... <w:p/>                     (same empty paragraph - untouched)
#       <w:p>                      (the empty paragraph + the "MVC" paragraph
#                                    merged into a single, bare paragraph)
#         <w:bookmarkStart w:id="0" w:name="_GoBack"/>
#         <w:bookmarkEnd w:id="0"/>
#       </w:p>

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$blankPara = $d.Paragraphs.Item($count - 1)
$mvcPara   = $d.Paragraphs.Item($count)

# Sanity-check we are looking at the right paragraphs before mutating.
if ($blankPara.Range.Text -ne "`r" -or $mvcPara.Range.Text -ne "MVC`r") {
    throw "Unexpected document tail - aborting so we don't corrupt content."
}

# 1) Merge the trailing blank paragraph into the "MVC" paragraph by
#    deleting the blank paragraph's own range (its paragraph mark).
$blankPara.Range.Delete()

# 2) The (now last) paragraph still holds "MVC" plus eastAsia rFonts
#    hints on both the run and the paragraph mark. Replace the whole
#    paragraph's contents with a clean paragraph that only keeps the
#    _GoBack bookmark, matching the target markup exactly.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
